# Multiply the wind power values in column D (rows 2-25) by 4,
# as described in the commit message: "4 times the wind power"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 4)   # column D
    $cell.Value = $cell.Value2 * 4
}
